# PlantillaLineasFacturaEmitida.xlsx edit
# - Rename column header "IdTipoDocumento" -> "TipoDocumento" (col A)
# - Rename column header "CodigoImpuesto" -> "CodigoEtax" (col T)
# - Give column T ("CodigoEtax") a custom width
# - Update the active selection / view state

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename header A1: IdTipoDocumento -> TipoDocumento
$ws.Range("A1").Value = "TipoDocumento"

# Rename header T1: CodigoImpuesto -> CodigoEtax
$ws.Range("T1").Value = "CodigoEtax"

# Give the newly (re)named CodigoEtax column (T) an explicit custom width,
# matching the width added for that column in the updated sheet.
$ws.Range("T1").ColumnWidth = 11.95

# Move the selection to match the workbook's saved view state.
$ws.Range("T13").Select()
